# Invalid Login - Data From XL
# Rename the existing sheet to "ValidLogin", add a new "InvalidLogin" sheet
# right after it, and populate both with UserName/Password/FailMsg data
# pulled from the external test data.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ValidLogin (the original sheet, renamed & populated) ---
$validLogin = $wb.ActiveSheet
$validLogin.Name = "ValidLogin"

$validLogin.Range("A1").Value = "UserName"
$validLogin.Range("B1").Value = "Password"
$validLogin.Range("C1").Value = "FailMsg"

$validLogin.Range("A2").Value = "admin"
$validLogin.Range("B2").Value = "manager"
$validLogin.Range("C2").Value = "Home Page is not displayed"

$validLogin.Columns.Item(1).AutoFit()
$validLogin.Columns.Item(3).AutoFit()

[void]$validLogin.Range("A1:C2").Select()

# --- Sheet 2: InvalidLogin (new sheet, added right after ValidLogin) ---
$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("C1").Value = "FailMsg"

$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"
$invalidLogin.Range("C2").Value = "Err MSg is not displayed"

$invalidLogin.Columns.Item(3).AutoFit()

[void]$invalidLogin.Range("C2").Select()
[void]$invalidLogin.Activate()
